$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): relabel B1 and add matching headers C1:E1 ---
# Copy B1's existing format (bold font, border, centered/top-aligned) onto the
# new header cells without disturbing the shared style table.
$ws.Range("B1").Copy()
$ws.Range("C1:E1").PasteSpecial(-4122)

$ws.Range("B1").Value2 = "YEAR"
$ws.Range("C1").Value2 = "Yards"
$ws.Range("D1").Value2 = "TDs"
$ws.Range("E1").Value2 = "GS"

# --- Row index column A: extend styled index values 7..10 for rows 9..12 ---
$ws.Range("A2").Copy()
$ws.Range("A9:A12").PasteSpecial(-4122)

$indexValues = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10)
for ($i = 0; $i -lt $indexValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $indexValues[$i]
}

# --- Data rows 2..12: YEAR / Yards / TDs / GS ---
$years  = @(2012, 2013, 2014, 2015, 2016, 2017, 2018, 2019, 2020, 2021, 2022)
$yards  = @(2316, 1911, 4862, 4000, 3812, 2039, 2588, 3691, 4449, 5476, 4593)
$tds    = @(11, 16, 27, 26, 22, 16, 14, 19, 25, 38, 27)
$gs     = @(37, 34, 63, 54, 52, 36, 31, 49, 55, 78, 56)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value2 = $years[$i]
    $ws.Cells.Item($row, 3).Value2 = $yards[$i]
    $ws.Cells.Item($row, 4).Value2 = $tds[$i]
    $ws.Cells.Item($row, 5).Value2 = $gs[$i]
}

Write-Output "done"
